# TestDataForPastFutureAbsences.xlsx was re-saved after tweaking the
# "DateQuestion" sample date on row 2 (used to exercise the Excel-read
# code path) and stripping the ad-hoc cell formatting / stale hyperlink
# that had accumulated on the fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the bespoke fonts/colors/date-format on every used cell - the
# fixture no longer needs the bold-header / blue-Courier styling.
$ws.Range("A1:I2").ClearFormats()

# Core data edit: push the DateQuestion sample value forward from
# 1/1/1992 to 11/11/1992 (serial 33604 -> 33919).
$ws.Range("B2").Value = "11/11/1992"
$ws.Range("B2").NumberFormat = "m/d/yyyy"

# The old mailto: hyperlink on the email cell is gone in the re-saved file.
$ws.Hyperlinks.Delete()

# Selection moved from I2 back to A2, and the sheet now carries explicit
# page-setup info instead.
$ws.Range("A2").Select()
$ws.PageSetup.Orientation = 1
